$d = $word.ActiveDocument
$wmain = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function New-PkgXml($bodyInner) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="' + $wmain + '"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) Rework the "IMPORTANTE" paragraph: merge the first two runs, and change
#    "no hay validación de datos" -> "de validacion no hay aumento de datos",
#    preserving spell-check proofErr wrapping around "dataset"/"validacion".
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("IMPORTANTE")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $newParaBody = (
        '<w:r><w:t xml:space="preserve">IMPORTANTE. En el </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>dataset</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> de </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>validacion</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> no hay </w:t></w:r>' +
        '<w:r><w:t>aumento</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> de datos. No se pueden </w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/><w:r><w:t>validador modelos</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t xml:space="preserve"> con datos que realmente no son de prueba</w:t></w:r>'
    )
    $target.Range.InsertXML((New-PkgXml ('<w:p>' + $newParaBody + '</w:p>'))) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Append a blank paragraph, a paragraph holding a page break, and another
#    blank paragraph at the very end of the document body.
# ---------------------------------------------------------------------------
$endPos = $d.Content.End
$insertPoint = $d.Range($endPos, $endPos)
$tailBody = '<w:p/><w:p><w:r><w:br w:type="page"/></w:r></w:p><w:p/>'
$insertPoint.InsertXML((New-PkgXml $tailBody)) | Out-Null
